# Round the "Maintenance Time" values in column B (the number inside the
# square brackets) to 1 decimal place, fixing the floating point noise
# that had accumulated (e.g. "[4.199999999999999]" -> "[4.2]").
# Cells whose value is already clean at 1 decimal place (e.g. "[4.0]",
# "[6]") are left untouched so their textual representation doesn't change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $text = $cell.Value2

    if ($text -ne $null -and $text -match '^\[(.+)\]$') {
        $inner = $Matches[1]
        $num = [double]$inner
        $rounded = [Math]::Round($num, 1)

        if ($rounded -ne $num) {
            $newText = "[" + $rounded + "]"
            $cell.Value = $newText
        }
    }
}
